$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge D2:D3
$ws.Range("D2:D3").UnMerge()

# Split the combined text into two cells
$ws.Range("D2").Value = "Ежедневно в 13:30 по рем.зоне."
$ws.Range("D2").VerticalAlignment = -4107

$ws.Range("D3").Value = "Ежедневно в 19:30 по магазину "
$ws.Range("D3").VerticalAlignment = -4107

# New merged placeholder area F2:F3 with vertical-center style (like old D2:D3 merge)
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4108
$ws.Range("F3").VerticalAlignment = -4108
$ws.Range("F2:F3").Merge()
